$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-03-30 Sunday" "2025-03-31 Monday"

Replace-Text "17×32=" "23×26="
Replace-Text "32×69=" "39×84="
Replace-Text "87×60=" "21×86="
Replace-Text "14×99=" "79×89="
Replace-Text "78×59=" "46×78="

Replace-Text "83×92=" "31×94="
Replace-Text "18×42=" "69×86="
Replace-Text "17×72=" "32×53="
Replace-Text "35×44=" "55×94="
Replace-Text "81×39=" "90×82="

Replace-Text "67×86=" "55×70="
Replace-Text "14×49=" "46×50="
Replace-Text "42×28=" "54×20="
Replace-Text "38×14=" "98×21="
Replace-Text "76×83=" "14×16="

Replace-Text "60×49=" "70×33="
Replace-Text "39×74=" "45×95="
Replace-Text "71×90=" "41×89="
Replace-Text "71×12=" "22×82="
Replace-Text "31×23=" "18×38="

Replace-Text "55×73=" "46×94="
Replace-Text "33×13=" "87×87="
Replace-Text "13×85=" "34×92="
Replace-Text "56×85=" "93×24="
Replace-Text "95×88=" "12×41="
